$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.006.95'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '2.470.21'
$ws.Range('E3').Value = '  +2.13%  '
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.70'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('D9').Value = '2.470.41'
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('E12').Value = '  +0.88%  '
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.96'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.80%  '
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').Value = '2.919.24'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('D17').Value = '62.963.57'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '2.473.09'
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.19'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '329.46'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('E22').Value = '  +8.98%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.23'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '664.57'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.92%  '
$ws.Range('E27').Value = '  +14.44%  '
$ws.Range('E28').Value = '  +0.37%  '
$ws.Range('D29').Value = '2.591.32'
$ws.Range('E30').Value = '  -9.31%  '
$ws.Range('E31').Value = '  +1.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.132'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.50%  '
$ws.Range('E35').Value = '  +3.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.77'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '152.55'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.371'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.76'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('E42').Value = '  -1.09%  '
$ws.Range('E43').Value = '  -0.57%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('E45').Value = '  +7.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '151.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.13'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +25.02%  '
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.65'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.93%  '
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0512'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.94%  '
